# Iterationsplan_as223jx.xlsx — add "Iteration 3" sheet, rename "Sheet3" -> "Iteration 2",
# and update the Iteration 2 sheet contents (per the commit "Lagt till Iterationsplan för
# Iteration 3 / Även lagt till mina Illustratorprojekt").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename "Sheet3" -> "Iteration 2"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet3")
$ws2.Name = "Iteration 2"

# ---------------------------------------------------------------------------
# 2. Update the "Iteration 2" sheet body to reflect newly-finished/changed work
#    (the two Illustrator-related items moved out, statuses progressed, and a
#    "verklig tid" (actual time) column E got filled in throughout).
# ---------------------------------------------------------------------------

# Row 14: "Gör en exempeldesign i Illustrator" -> renamed + marked Klar (done)
$ws2.Range("B14").Value = "Gör en exempeldesign för startsidan i Illustrator"
$ws2.Range("C14").Value = "Klar"
$ws2.Range("E14").Value = 2

# Row 15: "Designa ikon till appen" marked Klar (done)
$ws2.Range("C15").Value = "Klar"
$ws2.Range("E15").Value = 1

# Row 16: old "Spara-knappen" task replaced by the dropdown-list task, now Påbörjad
$ws2.Range("B16").Value = "Skapa dropdownlista för enheterna"
$ws2.Range("C16").Value = "Påbörjad"
$ws2.Range("E16").Value = 3

# Row 17: "Skapa Iterationsplan för Iteration 3" marked Klar (done)
$ws2.Range("C17").Value = "Klar"
$ws2.Range("E17").Value = 1

# Rows 12-13: fill in the new "verklig tid" (actual time) column E
$ws2.Range("E12").Value = 0
$ws2.Range("E13").Value = 0

# Summary rows 18-20: actual-time totals now live in column E instead of D
$ws2.Range("E18").Value = 12
$ws2.Range("E19").Value = 204
$ws2.Range("E20").Value = 240
$ws2.Range("D20").ClearContents()

# Column F widened now that it holds the "Kommentar" text for this sheet too
$ws2.Columns.Item(6).ColumnWidth = 19.67

# Selection moved off of B11 (old snapshot) onto F14
$ws2.Range("F14").Select()

# ---------------------------------------------------------------------------
# 3. Duplicate "Iteration 2" to create "Iteration 3" (keeps styles / merges /
#    column widths identical to start, matching how the author built the new
#    iteration sheet from the previous one), then rename + edit its content.
# ---------------------------------------------------------------------------
$ws2.Copy([System.Reflection.Missing]::Value, $ws2)
$ws3 = $wb.Worksheets.Item(4)
$ws3.Name = "Iteration 3"

# Row 3 text (comment about the just-finished iteration)
$ws3.Range("A3").Value = "Börjat jobba med min egen applikation. Insåg att dropdownlistan krånglade lite mer än jag förväntade mig, så det blir nog en del till tid för den."

# Row 5 text (goal for this iteration)
$ws3.Range("A5").Value = "Målet med denna Iteration är att jobba vidare med appen och förhoppningsvis bli klar med ""Skapa recept""-sidan."

# Table rows 8-15: new task list for Iteration 3, all "Ej påbörjad" (not started)
# and with no "verklig tid" column E recorded yet.
$ws3.Range("A8").ClearContents()
$ws3.Range("B8").Value = "Handledarmöte"
$ws3.Range("C8").Value = "Ej påbörjad"
$ws3.Range("D8").Value = 1
$ws3.Range("E8").ClearContents()

$ws3.Range("A9").Value = "F1"
$ws3.Range("B9").Value = "Göra klart dropdownen för enheter i Nytt Recept."
$ws3.Range("C9").Value = "Ej påbörjad"
$ws3.Range("D9").Value = 2
$ws3.Range("E9").ClearContents()

$ws3.Range("A10").Value = "F1/F2"
$ws3.Range("B10").Value = "Skicka vidare infon från textfälten till ny ""Visa Recept""-sida"
$ws3.Range("C10").Value = "Ej påbörjad"
$ws3.Range("D10").Value = 4
$ws3.Range("E10").ClearContents()

$ws3.Range("A11").Value = "F2"
$ws3.Range("B11").Value = "Strukturera upp ""Visa Recept""-sidan"
$ws3.Range("C11").Value = "Ej påbörjad"
$ws3.Range("D11").Value = 2
$ws3.Range("E11").ClearContents()

$ws3.Range("A12").Value = "F1"
$ws3.Range("B12").Value = "Gör exempeldesign för ""Nytt recept""-sidan i Illustrator"
$ws3.Range("C12").Value = "Ej påbörjad"
$ws3.Range("D12").Value = 2
$ws3.Range("E12").ClearContents()

$ws3.Range("A13").Value = "F1-F4"
$ws3.Range("B13").Value = "Läs ""Saving Data in SQL Databases"" på Androids developersida."
$ws3.Range("C13").Value = "Ej påbörjad"
$ws3.Range("D13").Value = 2
$ws3.Range("E13").ClearContents()

$ws3.Range("A14").ClearContents()
$ws3.Range("B14").Value = "Testrapport"
$ws3.Range("C14").Value = "Ej påbörjad"
$ws3.Range("D14").Value = 1
$ws3.Range("E14").ClearContents()

$ws3.Range("A15").ClearContents()
$ws3.Range("B15").Value = "Skapa Iterationsplan för Iteration 4"
$ws3.Range("C15").Value = "Ej påbörjad"
$ws3.Range("D15").Value = 1
$ws3.Range("E15").ClearContents()

# Row 16 ("Summa"): total estimated time rises to 15, no actual-time entry yet
$ws3.Range("D16").Value = 15
$ws3.Range("E16").ClearContents()

# Row 17 ("Tid sedan föregående iteration"): no value recorded yet
$ws3.Range("D17").ClearContents()
$ws3.Range("E17").ClearContents()

# Row 18 ("Total projekttid"): carried over from Iteration 2, now in column E
$ws3.Range("D18").ClearContents()
$ws3.Range("E18").Value = 240

# This sheet only has 18 used rows (vs. 20 on Iteration 2)
$ws3.Range("A19:F20").ClearContents()

# Column F is a bit narrower here since this sheet's comments are shorter
$ws3.Columns.Item(6).ColumnWidth = 15.74

# Selection lands on A15
$ws3.Range("A15").Select()

Write-Host "Edit complete"
